$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared string for the "P ... perinatal period" row (row 20) gets its
# trailing content extended with an extra "-" field. Writing the new value
# causes the (now-distinct) string to be appended as a new shared-string
# entry at the end, while every other row keeps referencing its original
# text - matching how Excel/the engine re-packs the shared-strings table.
$ws.Range("A20").Value = 'P                                         Certain conditions originating in the perinatal period            0                              1                                                                         -                                                                          - '

# Update the saved selection/active cell to A21 (as last set by the author).
$ws.Range("A21").Select()
